$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "J" column (spans 2:10) is no longer used anywhere else on the sheet;
# removing the stray J54 cell collapses rows 49-64 back down to spans "2:3".
$ws.Range("J54").ClearContents()

# Populate new rows 90-106 with the import/export admin-menu strings (EN/HE).
$ws.Range('B90').Value = 'Dashboard'
$ws.Range('C90').Value = 'אזור מנהל'
$ws.Range('B91').Value = 'All categories'
$ws.Range('C91').Value = 'כל המחלקות'
$ws.Range('B92').Value = 'Add New Category'
$ws.Range('C92').Value = 'הוסף מחלקה חדשה'
$ws.Range('B93').Value = 'Submit'
$ws.Range('C93').Value = 'אישור'
$ws.Range('B94').Value = ' Add Product'
$ws.Range('C94').Value = 'הוסף מוצר'
$ws.Range('B95').Value = ' Add New Product'
$ws.Range('C95').Value = 'הוסף מוצר חדש'
$ws.Range('B96').Value = 'Products'
$ws.Range('C96').Value = 'מוצרים'
$ws.Range('B97').Value = 'All Products'
$ws.Range('C97').Value = 'כל המוצרים'
$ws.Range('B98').Value = 'Add Product'
$ws.Range('C98').Value = 'הוסף מוצר'
$ws.Range('B99').Value = 'Deliveries'
$ws.Range('C99').Value = 'משלוחים'
$ws.Range('B100').Value = 'All Deliveries'
$ws.Range('C100').Value = 'כל המשלוחים'
$ws.Range('B101').Value = 'Add Delivery'
$ws.Range('C101').Value = 'הוסף משלוח'
$ws.Range('B102').Value = 'Orders'
$ws.Range('C102').Value = 'משלוחים'
$ws.Range('B103').Value = 'Customers'
$ws.Range('C103').Value = 'לקוחות'
$ws.Range('B104').Value = 'Configuration'
$ws.Range('C104').Value = 'תצורה'
$ws.Range('B105').Value = 'File import/export'
$ws.Range('C105').Value = 'ייבוא יצוא נתונים מקובץ'
$ws.Range('B106').Value = 'Search Product'
$ws.Range('C106').Value = 'חפש מוצר'

# Row 99 ("Deliveries") is a section header, like "Dashboard" in B4 - copy its
# header formatting onto B99 (format only, so the new value is left alone).
$ws.Range("B4").Copy()
$ws.Range("B99").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the sheet scrolled/selected on the newly added last row, matching the
# author-recorded selection.
$null = $ws.Range("B106").Select()
